# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" positioned between "总计" and "2022-Q2"
#    (built by duplicating "2022-Q2", which already has the right layout/
#    styles, then overwriting its data with the Q3 numbers).
# 2) Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q3 and push the existing 2022-Q2 / 2022-Q1 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q3" sheet just before "2022-Q2"
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2Index = $q2.Index
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item($q2Index)
$q3.Name = "2022-Q3"

# Row 2: 540002 / 汇丰晋信龙腾混合
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "540002"
$q3.Range("C2").Value = "汇丰晋信龙腾混合"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "4.72"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "93.98"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "4.31"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.2034"
$q3.Range("H2").Value = 10

# Row 3: 000892 / 九泰天宝灵活配置混合A
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "000892"
$q3.Range("C3").Value = "九泰天宝灵活配置混合A"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.06"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "89.35"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "7.85"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0047"
$q3.Range("H3").Value = 1

# Row 4 (brand new row): 002028 / 九泰天宝灵活配置混合C
# Copy A3's formatting (bold index-column style) down into A4 first.
$q3.Range("A3").Copy($q3.Range("A4"))
$q3.Range("A4").Value = 2
$q3.Range("B4").NumberFormat = "@"
$q3.Range("B4").Value = "002028"
$q3.Range("C4").Value = "九泰天宝灵活配置混合C"
$q3.Range("D4").NumberFormat = "@"
$q3.Range("D4").Value = "0.00"
$q3.Range("E4").NumberFormat = "@"
$q3.Range("E4").Value = "89.35"
$q3.Range("F4").NumberFormat = "@"
$q3.Range("F4").Value = "7.85"
$q3.Range("G4").Value = 0
$q3.Range("H4").Value = 1

# ---------------------------------------------------------------------
# Step 2: update "总计" with the new 2022-Q3 row on top, shifting the
# existing 2022-Q2 / 2022-Q1 rows down.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Extend the bold index-column style (currently on A2:A3) down to A4.
$zj.Range("A3").Copy($zj.Range("A4"))

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 3
$zj.Range("D2").Value = 0.21

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q2"
$zj.Range("C3").Value = 2
$zj.Range("D3").Value = 0.01

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2022-Q1"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.04

# Restore the originally active tab ("2022-Q1") since copying/renaming
# sheets above shifted focus onto the freshly created "2022-Q3" sheet.
$wb.Worksheets.Item("2022-Q1").Activate()
